$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21 ---
$ws.Range("B21").Value = 45376
$ws.Range("B21").NumberFormat = "m/d/yy"

$ws.Range("C21").Value = 0.41666666666666669
$ws.Range("C21").NumberFormat = "h:mm"

$ws.Range("D21").Value = 3

$ws.Range("G21").Value = "#Swapping Between Observer and Driver"
$ws.Range("E21").Value = "Observer/Driver"
$ws.Range("E21:F21").HorizontalAlignment = -4108
$ws.Range("E21:F21").MergeCells = $true

# --- Row 22 ---
$ws.Range("B22").Value = 45376
$ws.Range("B22").NumberFormat = "m/d/yy"

$ws.Range("C22").Value = 0.66666666666666663
$ws.Range("C22").NumberFormat = "h:mm"

$ws.Range("D22").Value = 2

$ws.Range("E22").Value = "Observer/Driver"
$ws.Range("E22:F22").HorizontalAlignment = -4108
$ws.Range("E22:F22").MergeCells = $true

$ws.Range("G22").Select()
